# Update countries & provincias Spain
# Applies the 24-May-2020 12:35 COVID data refresh to the "Pais" sheet:
#  - Updates the "last updated" timestamp in A1
#  - Refreshes case counts for a handful of countries (rows 30, 60, 63-65, 199-200)
#  - A few countries swapped rank with their neighbour, so the country label
#    (column A) and its stats move together to the row that now holds them

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / timestamp -----------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 12:35"

# --- Row 30: Suiza ------------------------------------------------------
$ws.Range("B30").Value = 30736
$ws.Range("C30").Value = 11
$ws.Range("E30").Value = 831

# --- Row 60: Marruecos ---------------------------------------------------
$ws.Range("B60").Value = 7429
$ws.Range("C60").Value = 23
$ws.Range("D60").Value = 4686
$ws.Range("E60").Value = 2545

# --- Row 63: Moldavia -----------------------------------------------------
$ws.Range("D63").Value = 3713
$ws.Range("E63").Value = 3031
$ws.Range("G63").Value = 8
$ws.Range("H63").Value = 250

# --- Row 64: was Armenia, now Ghana (rank swap) ---------------------------
$ws.Range("A64").Value = "Ghana"
$ws.Range("B64").Value = 6683
$ws.Range("C64").Value = 66
$ws.Range("D64").Value = 1998
$ws.Range("E64").Value = 4653
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 32

# --- Row 65: was Ghana, now Armenia (rank swap) ----------------------------
$ws.Range("A65").Value = "Armenia"
$ws.Range("B65").Value = 6661
$ws.Range("C65").Value = 359
$ws.Range("D65").Value = 3064
$ws.Range("E65").Value = 3516
$ws.Range("G65").Value = 4
$ws.Range("H65").Value = 81

# --- Row 199: was Belice, now Nueva Caledonia (rank swap) ------------------
$ws.Range("A199").Value = "Nueva Caledonia"
$ws.Range("D199").Value = 18
$ws.Range("H199").Value = 0

# --- Row 200: was Nueva Caledonia, now Belice (rank swap) ------------------
$ws.Range("A200").Value = "Belice"
$ws.Range("D200").Value = 16
$ws.Range("H200").Value = 2

# --- Rows 214-216: Sahara Occidental / San Bartolome / Bonaire shuffle -----
# (stats for these three rows are identical, only the labels rotate)
$ws.Range("A214").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A215").Value = "Sahara Occidental"
$ws.Range("A216").Value = "San Bartolome"
